$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill "Done" status in column E for rows 7 through 21 (previously blank)
$ws.Range("E7:E21").Value = "Done"

# Update the saved view state: scroll so row 4 is at top (A4), select G22
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G22").Select()
